# Remove every paragraph after the first "Bonjour, je m'appelle Peng"
# heading (Titre2 / numId 4) so the document starts its new "summary"
# structure with just that opening heading, directly followed by the
# section properties.
$d = $word.ActiveDocument

if ($d.Paragraphs.Count -gt 1) {
    $start = $d.Paragraphs(2).Range.Start
    $end = $d.Content.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
